# Update the NATMI ligand-receptor pair table (Proc-Procr) with refreshed TPM-derived
# expression / specificity values. Only the numeric result columns (E..T, excluding the
# constant K/L receptor-detection columns) need updating; identifiers in A-D are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending=ECs, Target=ECs
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.347467
$ws.Range("H2").Value = 1.042401
$ws.Range("I2").Value = 0.864291305025554
$ws.Range("J2").Value = 0.864291305025554
$ws.Range("M2").Value = 9.063968000000001
$ws.Range("N2").Value = 27.191904
$ws.Range("O2").Value = 0.3921806585023803
$ws.Range("P2").Value = 0.3921806585023803
$ws.Range("Q2").Value = 3.149429769056
$ws.Range("R2").Value = 28.344867921504
$ws.Range("S2").Value = 0.3389583331428034
$ws.Range("T2").Value = 0.3389583331428034

# Row 3: Sending=ECs, Target=Proc
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.347467
$ws.Range("H3").Value = 1.042401
$ws.Range("I3").Value = 0.864291305025554
$ws.Range("J3").Value = 0.864291305025554
$ws.Range("O3").Value = 0.185540687098555
$ws.Range("P3").Value = 0.185540687098555
$ws.Range("Q3").Value = 1.489995364765667
$ws.Range("R3").Value = 13.409958282891
$ws.Range("S3").Value = 0.160361202587748
$ws.Range("T3").Value = 0.1603612025877481

# Row 4: Sending=ECs, Target=MuSCs
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.347467
$ws.Range("H4").Value = 1.042401
$ws.Range("I4").Value = 0.864291305025554
$ws.Range("J4").Value = 0.864291305025554
$ws.Range("M4").Value = 9.518580666666667
$ws.Range("N4").Value = 28.555742
$ws.Range("O4").Value = 0.4118508840566691
$ws.Range("P4").Value = 0.4118508840566692
$ws.Range("Q4").Value = 3.307392668504666
$ws.Range("R4").Value = 29.766534016542
$ws.Range("S4").Value = 0.3559591380572667
$ws.Range("T4").Value = 0.3559591380572668

# Row 5: Sending=ECs, Target=FAPs
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.347467
$ws.Range("H5").Value = 1.042401
$ws.Range("I5").Value = 0.864291305025554
$ws.Range("J5").Value = 0.864291305025554
$ws.Range("M5").Value = 0.2410036666666667
$ws.Range("N5").Value = 0.7230110000000001
$ws.Range("O5").Value = 0.01042777034239546
$ws.Range("P5").Value = 0.01042777034239546
$ws.Range("Q5").Value = 0.08374082104566667
$ws.Range("R5").Value = 0.753667389411
$ws.Range("S5").Value = 0.009012631237735739
$ws.Range("T5").Value = 0.009012631237735741

# Row 6: Sending=MuSCs, Target=ECs
$ws.Range("I6").Value = 0.135708694974446
$ws.Range("J6").Value = 0.135708694974446
$ws.Range("M6").Value = 9.063968000000001
$ws.Range("N6").Value = 27.191904
$ws.Range("O6").Value = 0.3921806585023803
$ws.Range("P6").Value = 0.3921806585023803
$ws.Range("Q6").Value = 0.4945149874666667
$ws.Range("R6").Value = 4.4506348872
$ws.Range("S6").Value = 0.05322232535957692
$ws.Range("T6").Value = 0.05322232535957692

# Row 7: Sending=MuSCs, Target=Proc
$ws.Range("I7").Value = 0.135708694974446
$ws.Range("J7").Value = 0.135708694974446
$ws.Range("O7").Value = 0.185540687098555
$ws.Range("P7").Value = 0.185540687098555
$ws.Range("S7").Value = 0.02517948451080693
$ws.Range("T7").Value = 0.02517948451080693

# Row 8: Sending=MuSCs, Target=MuSCs
$ws.Range("I8").Value = 0.135708694974446
$ws.Range("J8").Value = 0.135708694974446
$ws.Range("M8").Value = 9.518580666666667
$ws.Range("N8").Value = 28.555742
$ws.Range("O8").Value = 0.4118508840566691
$ws.Range("P8").Value = 0.4118508840566692
$ws.Range("Q8").Value = 0.5193178968722222
$ws.Range("R8").Value = 4.67386107185
$ws.Range("S8").Value = 0.05589174599940246
$ws.Range("T8").Value = 0.05589174599940246

# Row 9: Sending=MuSCs, Target=FAPs
$ws.Range("I9").Value = 0.135708694974446
$ws.Range("J9").Value = 0.135708694974446
$ws.Range("M9").Value = 0.2410036666666667
$ws.Range("N9").Value = 0.7230110000000001
$ws.Range("O9").Value = 0.01042777034239546
$ws.Range("P9").Value = 0.01042777034239546
$ws.Range("Q9").Value = 0.01314875838055556
$ws.Range("R9").Value = 0.118338825425
$ws.Range("S9").Value = 0.00141513910465972
$ws.Range("T9").Value = 0.00141513910465972